# Trade #123 closed at 2026-02-17 16:07:04 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1199.62
$wsSummary.Range("B4").Value = -0.39
$wsSummary.Range("B5").Value = -0.06
$wsSummary.Range("B6").Value = 123
$wsSummary.Range("B7").Value = 47
$wsSummary.Range("B9").Value = 38.21

# --- Strategy Status sheet (MarketMaking row = row 4) ---
$wsStrategy = $wb.Worksheets.Item("Strategy Status")
$wsStrategy.Range("C4").Value = 99.62
$wsStrategy.Range("D4").Value = 123
$wsStrategy.Range("E4").Value = -0.39
$wsStrategy.Range("F4").Value = -0.38
$wsStrategy.Range("G4").Value = 38.21

# --- New trade row data (trade #123) ---
# Note: the Date column value looks like a date ("2026-02-17"), so it is
# prefixed with an apostrophe to force Excel to store it as plain text
# (matching the source data, which keeps it as a string, not a date serial).
$newRow = @(123, "'2026-02-17", "16:06:58", "MarketMaking", "UP", 0.3, 0.403551, "CLOSED", 34.5171, 0.1, 99.62, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.13)

# --- All Trades sheet ---
$wsAllTrades = $wb.Worksheets.Item("All Trades")
$rowAllTrades = 124
for ($i = 0; $i -lt $newRow.Length; $i++) {
    $wsAllTrades.Cells.Item($rowAllTrades, $i + 1).Value = $newRow[$i]
}

# --- MarketMaking sheet ---
$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
$rowMarketMaking = 124
for ($i = 0; $i -lt $newRow.Length; $i++) {
    $wsMarketMaking.Cells.Item($rowMarketMaking, $i + 1).Value = $newRow[$i]
}
